# Add an "image" column (E) of profile-picture URLs to the attendance sheet,
# with a few of the entries turned into clickable hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url1 = "https://img.freepik.com/free-photo/handsome-young-man-with-arms-crossed-white-background_23-2148222620.jpg"
$url2 = "https://img.freepik.com/free-photo/handsome-confident-smiling-man-with-hands-crossed-chest_176420-18743.jpg"
$url3 = "https://img.freepik.com/free-photo/attractive-mixed-race-male-with-positive-smile-shows-white-teeth-keeps-hands-stomach-being-high-spirit-wears-white-shirt-rejoices-positive-moments-life-people-emotions-concept_273609-15527.jpg"
$url4 = "https://img.freepik.com/free-photo/confident-handsome-guy-posing-against-white-wall_176420-32936.jpg"
$url5 = "https://img.freepik.com/free-photo/fashionable-stylish-man-with-dark-eyes-casual-clothes-looking-aside-with-placid-thoughtful-look-pensive-guy-with-puzzled-expression-thinking-about-something-building-plans_176420-10331.jpg"
$url6 = "https://img.freepik.com/free-photo/thoughtful-concerned-man-thinking-trying-find-solution_176420-19574.jpg"
$url7 = "https://img.freepik.com/free-photo/serious-thoughtful-man-squinting-skeptical-thinking-as-making-choice_176420-19020.jpg"

# Header
$ws.Range("E1").Value = "image"

# One URL per attendee (rows 2-11); some values repeat further down the list
$ws.Range("E2").Value  = $url1
$ws.Range("E3").Value  = $url2
$ws.Range("E4").Value  = $url3
$ws.Range("E5").Value  = $url4
$ws.Range("E6").Value  = $url5
$ws.Range("E7").Value  = $url6
$ws.Range("E8").Value  = $url7
$ws.Range("E9").Value  = $url3
$ws.Range("E10").Value = $url6
$ws.Range("E11").Value = $url1

# Turn a few of the image cells into live hyperlinks (pointing at the same URL
# shown as text), mirroring the state of the edited workbook.
$ws.Hyperlinks.Add($ws.Range("E4"), $url3)
$ws.Hyperlinks.Add($ws.Range("E7"), $url6)
$ws.Hyperlinks.Add($ws.Range("E2"), $url1)

# The sheet was left with the new column selected.
$ws.Range("E1:E11").Select()
